# Auto-generated script applying the Anima_Profits.xlsx diff
# Updates currentAveragePrice / LevePrice / LeveProfit columns (H-N) across all 8 job sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 920.0333000000001
$ws.Range("I6").Value = 104.04
$ws.Range("J6").Value = 5000
$ws.Range("K6").Value = 312.12
$ws.Range("L6").Value = 15000
$ws.Range("M6").Value = -200.12
$ws.Range("N6").Value = -15224
$ws.Range("H8").Value = 293.16666
$ws.Range("I8").Value = 293.16666
$ws.Range("K8").Value = 879.4999799999999
$ws.Range("M8").Value = -740.4999799999999
$ws.Range("H112").Value = 11376.0625
$ws.Range("I112").Value = 0
$ws.Range("J112").Value = 11376.0625
$ws.Range("K112").Value = 0
$ws.Range("L112").Value = 34128.1875
$ws.Range("M112").ClearContents()
$ws.Range("N112").Value = -36344.1875
$ws.Range("H116").Value = 3039.5833
$ws.Range("I116").Value = 2780.5
$ws.Range("J116").Value = 4335
$ws.Range("K116").Value = 2780.5
$ws.Range("L116").Value = 4335
$ws.Range("M116").Value = 661.5
$ws.Range("N116").Value = -11219
$ws.Range("H138").Value = 1899.2693
$ws.Range("I138").Value = 1766.6
$ws.Range("J138").Value = 2022.1111
$ws.Range("K138").Value = 5299.799999999999
$ws.Range("L138").Value = 6066.3333
$ws.Range("M138").Value = -159.7999999999993
$ws.Range("N138").Value = -16346.3333
$ws.Range("H141").Value = 4364.9395
$ws.Range("I141").Value = 1171.9131
$ws.Range("J141").Value = 11708.9
$ws.Range("K141").Value = 3515.7393
$ws.Range("L141").Value = 35126.7
$ws.Range("M141").Value = 1664.2607
$ws.Range("N141").Value = -45486.7

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 471522.7
$ws.Range("I32").Value = 514067.94
$ws.Range("K32").Value = 514067.94
$ws.Range("M32").Value = -513780.94
$ws.Range("H45").Value = 2662.6667
$ws.Range("I45").Value = 2733.6
$ws.Range("J45").Value = 2544.4443
$ws.Range("K45").Value = 2733.6
$ws.Range("L45").Value = 2544.4443
$ws.Range("M45").Value = -2356.6
$ws.Range("N45").Value = -3298.4443
$ws.Range("H74").Value = 1268.8667
$ws.Range("I74").Value = 968.3333
$ws.Range("J74").Value = 2471
$ws.Range("K74").Value = 968.3333
$ws.Range("L74").Value = 2471
$ws.Range("M74").Value = -94.33330000000001
$ws.Range("N74").Value = -4219
$ws.Range("H77").Value = 1268.8667
$ws.Range("I77").Value = 968.3333
$ws.Range("J77").Value = 2471
$ws.Range("K77").Value = 4841.6665
$ws.Range("L77").Value = 12355
$ws.Range("M77").Value = -473.6665000000003
$ws.Range("N77").Value = -21091
$ws.Range("H92").Value = 73750
$ws.Range("J92").Value = 73750
$ws.Range("L92").Value = 73750
$ws.Range("N92").Value = -78742
$ws.Range("H101").Value = 79602
$ws.Range("J101").Value = 79602
$ws.Range("L101").Value = 79602
$ws.Range("N101").Value = -86092
$ws.Range("H102").Value = 2233.3635
$ws.Range("I102").Value = 2233.3635
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 2233.3635
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = -611.3634999999999
$ws.Range("N102").ClearContents()
$ws.Range("H105").Value = 60000
$ws.Range("J105").Value = 60000
$ws.Range("L105").Value = 60000
$ws.Range("N105").Value = -66988
$ws.Range("H119").Value = 39800
$ws.Range("J119").Value = 39800
$ws.Range("L119").Value = 39800
$ws.Range("N119").Value = -49476
$ws.Range("H122").Value = 127454.25
$ws.Range("I122").Value = 251953
$ws.Range("J122").Value = 2955.5
$ws.Range("K122").Value = 755859
$ws.Range("L122").Value = 8866.5
$ws.Range("M122").Value = -753409
$ws.Range("N122").Value = -13766.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H7").Value = 498
$ws.Range("I7").Value = 498
$ws.Range("K7").Value = 498
$ws.Range("M7").Value = -385
$ws.Range("H134").Value = 3050.35
$ws.Range("I134").Value = 2962.4167
$ws.Range("J134").Value = 3182.25
$ws.Range("K134").Value = 8887.250100000001
$ws.Range("L134").Value = 9546.75
$ws.Range("M134").Value = -6352.250100000001
$ws.Range("N134").Value = -14616.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 426.66666
$ws.Range("I22").Value = 395.83334
$ws.Range("J22").Value = 550
$ws.Range("K22").Value = 395.83334
$ws.Range("L22").Value = 550
$ws.Range("M22").Value = -45.83334000000002
$ws.Range("N22").Value = -1250
$ws.Range("H99").Value = 1745.1852
$ws.Range("I99").Value = 1324.2858
$ws.Range("J99").Value = 1892.5
$ws.Range("K99").Value = 1324.2858
$ws.Range("L99").Value = 1892.5
$ws.Range("M99").Value = 173.7141999999999
$ws.Range("N99").Value = -4888.5
$ws.Range("H126").Value = 1745.1852
$ws.Range("I126").Value = 1324.2858
$ws.Range("J126").Value = 1892.5
$ws.Range("K126").Value = 3972.8574
$ws.Range("L126").Value = 5677.5
$ws.Range("M126").Value = -1502.8574
$ws.Range("N126").Value = -10617.5
$ws.Range("H132").Value = 9806442
$ws.Range("I132").Value = 2179.8462
$ws.Range("K132").Value = 6539.5386
$ws.Range("M132").Value = -4009.5386
$ws.Range("H134").Value = 2833.75
$ws.Range("I134").Value = 1343.5
$ws.Range("J134").Value = 5317.5
$ws.Range("K134").Value = 4030.5
$ws.Range("L134").Value = 15952.5
$ws.Range("M134").Value = -1495.5
$ws.Range("N134").Value = -21022.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 452.95
$ws.Range("J5").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("N5").ClearContents()
$ws.Range("H104").Value = 3938.1667
$ws.Range("I104").Value = 2800
$ws.Range("K104").Value = 8400
$ws.Range("M104").Value = -5779
$ws.Range("H121").Value = 358183.4
$ws.Range("I121").Value = 1667018
$ws.Range("J121").Value = 1228.5
$ws.Range("K121").Value = 5001054
$ws.Range("L121").Value = 3685.5
$ws.Range("M121").Value = -4999744
$ws.Range("N121").Value = -6305.5
$ws.Range("H122").Value = 9759.272000000001
$ws.Range("I122").Value = 594.8889
$ws.Range("J122").Value = 50999
$ws.Range("K122").Value = 5354.0001
$ws.Range("L122").Value = 458991
$ws.Range("M122").Value = -2904.0001
$ws.Range("N122").Value = -463891
$ws.Range("H131").Value = 1021.9677
$ws.Range("I131").Value = 416.66666
$ws.Range("J131").Value = 1086.8214
$ws.Range("K131").Value = 1249.99998
$ws.Range("L131").Value = 3260.4642
$ws.Range("M131").Value = 3790.00002
$ws.Range("N131").Value = -13340.4642
$ws.Range("H135").Value = 452.95
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()
$ws.Range("H137").Value = 8342680.5
$ws.Range("I137").Value = 20848538
$ws.Range("J137").Value = 5442.5835
$ws.Range("K137").Value = 62545614
$ws.Range("L137").Value = 16327.7505
$ws.Range("M137").Value = -62540514
$ws.Range("N137").Value = -26527.7505

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 46311804
$ws.Range("I80").Value = 63627980
$ws.Range("J80").Value = 135333.33
$ws.Range("K80").Value = 63627980
$ws.Range("L80").Value = 135333.33
$ws.Range("M80").Value = -63626982
$ws.Range("N80").Value = -137329.33
$ws.Range("H83").Value = 46311804
$ws.Range("I83").Value = 63627980
$ws.Range("J83").Value = 135333.33
$ws.Range("K83").Value = 318139900
$ws.Range("L83").Value = 676666.6499999999
$ws.Range("M83").Value = -318134908
$ws.Range("N83").Value = -686650.6499999999
$ws.Range("H102").Value = 1319.1428
$ws.Range("I102").Value = 1199.2142
$ws.Range("J102").Value = 1559
$ws.Range("K102").Value = 1199.2142
$ws.Range("L102").Value = 1559
$ws.Range("M102").Value = 422.7858000000001
$ws.Range("N102").Value = -4803
$ws.Range("H122").Value = 2550
$ws.Range("I122").Value = 2000
$ws.Range("J122").Value = 2733.3333
$ws.Range("K122").Value = 6000
$ws.Range("L122").Value = 8199.999899999999
$ws.Range("M122").Value = -3550
$ws.Range("N122").Value = -13099.9999
$ws.Range("H132").Value = 3359.6667
$ws.Range("I132").Value = 3143.2727
$ws.Range("J132").Value = 3699.7144
$ws.Range("K132").Value = 9429.8181
$ws.Range("L132").Value = 11099.1432
$ws.Range("M132").Value = -6899.8181
$ws.Range("N132").Value = -16159.1432

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 3808.5557
$ws.Range("I132").Value = 3371.6
$ws.Range("J132").Value = 4354.75
$ws.Range("K132").Value = 10114.8
$ws.Range("L132").Value = 13064.25
$ws.Range("M132").Value = -7584.799999999999
$ws.Range("N132").Value = -18124.25
$ws.Range("H136").Value = 2378.4546
$ws.Range("I136").Value = 2205.1428
$ws.Range("J136").Value = 2681.75
$ws.Range("K136").Value = 6615.428400000001
$ws.Range("L136").Value = 8045.25
$ws.Range("M136").Value = -4065.428400000001
$ws.Range("N136").Value = -13145.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H92").Value = 38333.332
$ws.Range("J92").Value = 38333.332
$ws.Range("L92").Value = 38333.332
$ws.Range("N92").Value = -43325.332
$ws.Range("H123").Value = 22699
$ws.Range("J123").Value = 23238.8
$ws.Range("L123").Value = 23238.8
$ws.Range("N123").Value = -33038.8
$ws.Range("H136").Value = 2109.3171
$ws.Range("I136").Value = 1607.0714
$ws.Range("J136").Value = 3191.077
$ws.Range("K136").Value = 4821.2142
$ws.Range("L136").Value = 9573.231
$ws.Range("M136").Value = -2271.2142
$ws.Range("N136").Value = -14673.231

Write-Host "Applied Anima_Profits.xlsx updates to all sheets"